$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.978.37"
$ws.Range("E2").Value = "  -1.85%  "
$ws.Range("D3").Value = "3.806.95"
$ws.Range("E3").Value = "  +3.28%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "621.65"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.15"
$ws.Range("E6").Value = "  -3.80%  "
$ws.Range("D7").Value = "3.804.71"
$ws.Range("E7").Value = "  +3.32%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.536"
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.172"
$ws.Range("E10").Value = "  +4.47%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.29"
$ws.Range("E11").Value = "  -4.97%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.493"
$ws.Range("E12").Value = "  -1.69%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.69"
$ws.Range("E13").Value = "  +1.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000263"
$ws.Range("E14").Value = "  +3.09%  "
$ws.Range("D15").Value = "4.440.79"
$ws.Range("E15").Value = "  +3.21%  "
$ws.Range("D16").Value = "3.803.52"
$ws.Range("E16").Value = "  +3.17%  "
$ws.Range("D17").Value = "70.040.39"
$ws.Range("E17").Value = "  -1.76%  "
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.59"
$ws.Range("E19").Value = "  +1.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.87"
$ws.Range("E20").Value = "  -0.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "510.06"
$ws.Range("E21").Value = "  -1.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.58"
$ws.Range("E22").Value = "  +4.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.728"
$ws.Range("E23").Value = "  -2.50%  "
$ws.Range("B24").Value = "Fetch.AI"
$ws.Range("C24").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.51"
$ws.Range("E24").Value = "  +3.07%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "87.78"
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.25"
$ws.Range("E26").Value = "  -1.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000142"
$ws.Range("E27").Value = "  +29.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.05"
$ws.Range("E28").Value = "  +1.30%  "
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.50"
$ws.Range("E30").Value = "  -1.60%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.88"
$ws.Range("E31").Value = "  +3.99%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.77"
$ws.Range("E32").Value = "  -5.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.53"
$ws.Range("E33").Value = "  -0.53%  "
$ws.Range("E34").Value = "  -1.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.08"
$ws.Range("E36").Value = "  +6.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.22"
$ws.Range("E37").Value = "  +1.13%  "
$ws.Range("E38").Value = "  +4.92%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.332"
$ws.Range("E39").Value = "  -3.26%  "
$ws.Range("E40").Value = "  -1.38%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "51.07"
$ws.Range("E41").Value = "  +0.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "45.42"
$ws.Range("E42").Value = "  -0.98%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.74"
$ws.Range("E43").Value = "  -1.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "417.43"
$ws.Range("E44").Value = "  +2.53%  "
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.83"
$ws.Range("E45").Value = "  +1.49%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "3.041.75"
$ws.Range("E46").Value = "  -4.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0363"
$ws.Range("E47").Value = "  -1.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.38"
$ws.Range("E48").Value = "  -2.84%  "
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "137.84"
$ws.Range("E50").Value = "  +0.81%  "
$ws.Range("E51").Value = "  +0.82%  "
